# Mapear pessoas por turnos - Garantir 2 pessoas de cada departamento sempre
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extra rows: shift-number, Equipa, Disponibilidade (duplicating the
# existing department/availability pairs from rows 2-10, but keyed by
# a running shift number instead of a person's name).
$data = @(
    @(1, "Webdev",   "[1,1,1,1,1,0,1,1,1,0]"),
    @(2, "Webdev",   "[1,1,0,1,1,1,1,1,1,1]"),
    @(3, "Webdev",   "[1,1,1,1,0,1,1,1,0,1]"),
    @(4, "Speakers", "[1,1,0,1,0,0,1,1,1,1]"),
    @(5, "Webdev",   "[1,1,0,1,1,1,1,1,1,1]"),
    @(6, "Speakers", "[1,1,0,1,1,1,1,1,1,1]"),
    @(7, "Buss",     "[0,1,1,1,1,1,1,1,1,1]"),
    @(8, "Webdev",   "[1,1,0,1,1,1,1,1,1,1]"),
    @(9, "Buss",     "[1,1,1,1,1,1,1,1,1,1]")
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$ws.Range("A20").Select()
